# update HKStock add function readIndexData
# Appends three new trading-day rows (index 263-265, sheet rows 265-267) to Sheet1,
# mirroring the shape/style of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @{ Row = 265; Idx = 263; Date = "2016-08-09"; Open = 13427.37; Low = 13419.69; High = 13477.47; Close = 13459.27 },
    @{ Row = 266; Idx = 264; Date = "2016-08-10"; Open = 13453.64; Low = 13395.19; High = 13532.5;  Close = 13456.52 },
    @{ Row = 267; Idx = 265; Date = "2016-08-11"; Open = 13407.52; Low = 13367.43; High = 13530.9;  Close = 13453.74 }
)

# Copy the formatting (bold/border/center) of the last existing index cell
# down onto the new index cells so the new rows render just like the old ones.
$ws.Range("A264").Copy()
$ws.Range("A265:A267").PasteSpecial(-4122)  # xlPasteFormats

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value2 = $r.Idx
    $ws.Cells.Item($r.Row, 2).Value2 = "HSC"

    # Write the trade-date as a formula that evaluates to a text string so
    # Excel does not auto-convert the "yyyy-mm-dd" literal into a date
    # serial number, then paste only the resulting value back as plain text
    # (this keeps the cell's style untouched, matching the existing rows).
    $ws.Range("ZZ1").Formula = '="' + $r.Date + '"'
    $ws.Range("ZZ1").Copy()
    $ws.Cells.Item($r.Row, 3).PasteSpecial(-4163)  # xlPasteValues

    $ws.Cells.Item($r.Row, 4).Value2 = $r.Open
    $ws.Cells.Item($r.Row, 5).Value2 = $r.Low
    $ws.Cells.Item($r.Row, 6).Value2 = $r.High
    $ws.Cells.Item($r.Row, 7).Value2 = $r.Close
}

$ws.Range("ZZ1").ClearContents()
$excel.CutCopyMode = 0
